$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.0275626423690205
$ws.Range("J2").Value = 0.02947285804177352
$ws.Range("M2").Value = 6.780879499999999
$ws.Range("N2").Value = 13.561759
$ws.Range("O2").Value = 0.09314755032665376
$ws.Range("P2").Value = 0.07273600820493056
$ws.Range("Q2").Value = 0.009845837034
$ws.Range("R2").Value = 0.059075022204
$ws.Range("S2").Value = 0.002567392617203896
$ws.Range("T2").Value = 0.002143738044349192
$ws.Range("I3").Value = 0.0275626423690205
$ws.Range("J3").Value = 0.02947285804177352
$ws.Range("O3").Value = 0.3322252662272683
$ws.Range("P3").Value = 0.389136476570504
$ws.Range("S3").Value = 0.009157006198974822
$ws.Range("T3").Value = 0.01146896413283839
$ws.Range("I4").Value = 0.0275626423690205
$ws.Range("J4").Value = 0.02947285804177352
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.592065666666667
$ws.Range("N4").Value = 4.776197
$ws.Range("O4").Value = 0.02186987938794236
$ws.Range("P4").Value = 0.02561625701948875
$ws.Range("Q4").Value = 0.002311679348
$ws.Range("R4").Value = 0.020805114132
$ws.Range("S4").Value = 0.0006027916642234684
$ws.Range("T4").Value = 0.0007549843066969764
$ws.Range("I5").Value = 0.0275626423690205
$ws.Range("J5").Value = 0.02947285804177352
$ws.Range("M5").Value = 25.158886
$ws.Range("N5").Value = 50.31777200000001
$ws.Range("O5").Value = 0.3456024546443489
$ws.Range("P5").Value = 0.2698701456828592
$ws.Range("Q5").Value = 0.03653070247200001
$ws.Range("R5").Value = 0.219184214832
$ws.Range("S5").Value = 0.009525716859217818
$ws.Range("T5").Value = 0.007953844493423646
$ws.Range("I6").Value = 0.0275626423690205
$ws.Range("J6").Value = 0.02947285804177352
$ws.Range("M6").Value = 14.00046133333333
$ws.Range("N6").Value = 42.001384
$ws.Range("O6").Value = 0.1923214645892228
$ws.Range("P6").Value = 0.225266723235713
$ws.Range("Q6").Value = 0.020328669856
$ws.Range("R6").Value = 0.182958028704
$ws.Range("S6").Value = 0.005300887748358989
$ws.Range("T6").Value = 0.006639254155461653
$ws.Range("I7").Value = 0.0275626423690205
$ws.Range("J7").Value = 0.02947285804177352
$ws.Range("M7").Value = 1.079828666666667
$ws.Range("N7").Value = 3.239486
$ws.Range("O7").Value = 0.01483338482456395
$ws.Range("P7").Value = 0.01737438928650463
$ws.Range("Q7").Value = 0.001567911224
$ws.Range("R7").Value = 0.014111201016
$ws.Range("S7").Value = 0.000408847281041512
$ws.Range("T7").Value = 0.0005120729090036616
$ws.Range("I8").Value = 0.1944381169324222
$ws.Range("J8").Value = 0.1386090380724913
$ws.Range("M8").Value = 6.780879499999999
$ws.Range("N8").Value = 13.561759
$ws.Range("O8").Value = 0.09314755032665376
$ws.Range("P8").Value = 0.07273600820493056
$ws.Range("Q8").Value = 0.06945654871849999
$ws.Range("R8").Value = 0.277826194874
$ws.Range("S8").Value = 0.01811143428238258
$ws.Range("T8").Value = 0.01008186813051826
$ws.Range("I9").Value = 0.1944381169324222
$ws.Range("J9").Value = 0.1386090380724913
$ws.Range("O9").Value = 0.3322252662272683
$ws.Range("P9").Value = 0.389136476570504
$ws.Range("S9").Value = 0.06459725516260267
$ws.Range("T9").Value = 0.05393783269635611
$ws.Range("I10").Value = 0.1944381169324222
$ws.Range("J10").Value = 0.1386090380724913
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.592065666666667
$ws.Range("N10").Value = 4.776197
$ws.Range("O10").Value = 0.02186987938794236
$ws.Range("P10").Value = 0.02561625701948875
$ws.Range("Q10").Value = 0.01630752862366667
$ws.Range("R10").Value = 0.097845171742
$ws.Range("S10").Value = 0.004252338165730707
$ws.Range("T10").Value = 0.003550644744489039
$ws.Range("I11").Value = 0.1944381169324222
$ws.Range("J11").Value = 0.1386090380724913
$ws.Range("M11").Value = 25.158886
$ws.Range("N11").Value = 50.31777200000001
$ws.Range("O11").Value = 0.3456024546443489
$ws.Range("P11").Value = 0.2698701456828592
$ws.Range("Q11").Value = 0.2577024692980001
$ws.Range("R11").Value = 1.030809877192
$ws.Range("S11").Value = 0.06719829048827004
$ws.Range("T11").Value = 0.0374064412975842
$ws.Range("I12").Value = 0.1944381169324222
$ws.Range("J12").Value = 0.1386090380724913
$ws.Range("M12").Value = 14.00046133333333
$ws.Range("N12").Value = 42.001384
$ws.Range("O12").Value = 0.1923214645892228
$ws.Range("P12").Value = 0.225266723235713
$ws.Range("Q12").Value = 0.1434067254373333
$ws.Range("R12").Value = 0.8604403526240001
$ws.Range("S12").Value = 0.03739462342041399
$ws.Range("T12").Value = 0.03122400381744431
$ws.Range("I13").Value = 0.1944381169324222
$ws.Range("J13").Value = 0.1386090380724913
$ws.Range("M13").Value = 1.079828666666667
$ws.Range("N13").Value = 3.239486
$ws.Range("O13").Value = 0.01483338482456395
$ws.Range("P13").Value = 0.01737438928650463
$ws.Range("Q13").Value = 0.01106068503266667
$ws.Range("R13").Value = 0.066364110196
$ws.Range("S13").Value = 0.002884175413022181
$ws.Range("T13").Value = 0.002408247386099405
$ws.Range("G14").Value = 0.040985
$ws.Range("H14").Value = 0.122955
$ws.Range("I14").Value = 0.7779992406985573
$ws.Range("J14").Value = 0.8319181038857351
$ws.Range("M14").Value = 6.780879499999999
$ws.Range("N14").Value = 13.561759
$ws.Range("O14").Value = 0.09314755032665376
$ws.Range("P14").Value = 0.07273600820493056
$ws.Range("Q14").Value = 0.2779143463074999
$ws.Range("R14").Value = 1.667486077845
$ws.Range("S14").Value = 0.07246872342706727
$ws.Range("T14").Value = 0.0605104020300631
$ws.Range("G15").Value = 0.040985
$ws.Range("H15").Value = 0.122955
$ws.Range("I15").Value = 0.7779992406985573
$ws.Range("J15").Value = 0.8319181038857351
$ws.Range("O15").Value = 0.3322252662272683
$ws.Range("P15").Value = 0.389136476570504
$ws.Range("Q15").Value = 0.9912248617016667
$ws.Range("R15").Value = 8.921023755315
$ws.Range("S15").Value = 0.2584710048656908
$ws.Range("T15").Value = 0.3237296797413095
$ws.Range("G16").Value = 0.040985
$ws.Range("H16").Value = 0.122955
$ws.Range("I16").Value = 0.7779992406985573
$ws.Range("J16").Value = 0.8319181038857351
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.592065666666667
$ws.Range("N16").Value = 4.776197
$ws.Range("O16").Value = 0.02186987938794236
$ws.Range("P16").Value = 0.02561625701948875
$ws.Range("Q16").Value = 0.06525081134833333
$ws.Range("R16").Value = 0.5872573021349999
$ws.Range("S16").Value = 0.01701474955798819
$ws.Range("T16").Value = 0.02131062796830274
$ws.Range("G17").Value = 0.040985
$ws.Range("H17").Value = 0.122955
$ws.Range("I17").Value = 0.7779992406985573
$ws.Range("J17").Value = 0.8319181038857351
$ws.Range("M17").Value = 25.158886
$ws.Range("N17").Value = 50.31777200000001
$ws.Range("O17").Value = 0.3456024546443489
$ws.Range("P17").Value = 0.2698701456828592
$ws.Range("Q17").Value = 1.03113694271
$ws.Range("R17").Value = 6.18682165626
$ws.Range("S17").Value = 0.2688784472968611
$ws.Range("T17").Value = 0.2245098598918513
$ws.Range("G18").Value = 0.040985
$ws.Range("H18").Value = 0.122955
$ws.Range("I18").Value = 0.7779992406985573
$ws.Range("J18").Value = 0.8319181038857351
$ws.Range("M18").Value = 14.00046133333333
$ws.Range("N18").Value = 42.001384
$ws.Range("O18").Value = 0.1923214645892228
$ws.Range("P18").Value = 0.225266723235713
$ws.Range("Q18").Value = 0.5738089077466667
$ws.Range("R18").Value = 5.16428016972
$ws.Range("S18").Value = 0.1496259534204498
$ws.Range("T18").Value = 0.187403465262807
$ws.Range("G19").Value = 0.040985
$ws.Range("H19").Value = 0.122955
$ws.Range("I19").Value = 0.7779992406985573
$ws.Range("J19").Value = 0.8319181038857351
$ws.Range("M19").Value = 1.079828666666667
$ws.Range("N19").Value = 3.239486
$ws.Range("O19").Value = 0.01483338482456395
$ws.Range("P19").Value = 0.01737438928650463
$ws.Range("Q19").Value = 0.04425677790333333
$ws.Range("R19").Value = 0.39831100113
$ws.Range("S19").Value = 0.01154036213050025

# Sending cluster for rows 14-19 is re-labelled from "Resolving-Mac" to
# "Neutrophils" (the shared-string table slot that these rows pointed at was
# repurposed for the new cluster name; the net, user-visible effect is that
# the sending cluster column changes even though the cells own XML element
# was untouched in the source diff).
$ws.Range("A14").Value = "Neutrophils"
$ws.Range("A15").Value = "Neutrophils"
$ws.Range("A16").Value = "Neutrophils"
$ws.Range("A17").Value = "Neutrophils"
$ws.Range("A18").Value = "Neutrophils"
$ws.Range("A19").Value = "Neutrophils"
